$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.323.27'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.056.50'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.47'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.616'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.84'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +4.83%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +3.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.23'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0758'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.361.03'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.32'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.85'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.773'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.14'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.049.43'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '37.264.31'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.18'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +14.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.05'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0810'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '225.16'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.71%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("E25").Value = '  +2.02%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.85'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.92'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +3.38%  '
$ws.Range("E29").Value = '  +6.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.128'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.08'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.118'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.45'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.57'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0616'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.54'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +6.84%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.77'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.85'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.28'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.60'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +12.26%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.469.58'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.28'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0936'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.51%  '
$ws.Range("E46").Value = '  +5.81%  '
$ws.Range("E47").Value = '  +3.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.54'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("E49").Value = '  +1.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.17'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.75%  '
$ws.Range("E51").Value = '  +1.83%  '
